{"js": "// Apply the five textual edits described by the diff using the Word\n// JavaScript API. Each edit is performed as an independent\n// search-and-replace against context.document.body so the script is\n// resilient to however the original text happens to be split across\n// runs in the underlying OOXML.\n\nconst edits = [\n  {\n    find: \"IMPACT OF AUTO-GRADING OF AN INTRODUCTORY COMPUTING COURSE\",\n    replace: \"IMPACT OF AUTO-GRADING ON AN INTRODUCTORY COMPUTING COURSE\",\n  },\n  {\n    find:\n      \"make more submissions per assignment, and have qualitative evidence of students leveraging feedback to improve their programs.\",\n    replace:\n      \"make more submissions per assignment, indicating that students were leveraging feedback to improve their programs.\",\n  },\n  {\n    find:\n      \" (2) small case studies, where a sample of students\\u2019 individual paths through particular assignments were qualitatively observed and described.\",\n    replace:\n      \" (2) fine-grain submission rate analysis, where the distribution of submission rates among students was visible.\",\n  },\n  {\n    find: \"The aggregate analysis used data from seven offerings\",\n    replace: \"The analysis used data from seven offerings\",\n  },\n  {\n    find: \"The lack of any higher average submission rates indicates\",\n    replace: \"The lack of significantly higher average submission rates indicates\",\n  },\n];\n\nfor (const { find, replace } of edits) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the five textual edits described by the diff using the Word\n# COM object model. Each edit is performed as an independent\n# Find/Replace against the whole document so the script does not\n# depend on exactly how the original text happens to be split across\n# runs in the underlying OOXML.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceOne)\n    )\n}\n\nReplace-Text \"IMPACT OF AUTO-GRADING OF AN INTRODUCTORY COMPUTING COURSE\" \"IMPACT OF AUTO-GRADING ON AN INTRODUCTORY COMPUTING COURSE\"\n\nReplace-Text \"make more submissions per assignment, and have qualitative evidence of students leveraging feedback to improve their programs.\" \"make more submissions per assignment, indicating that students were leveraging feedback to improve their programs.\"\n\nReplace-Text \" (2) small case studies, where a sample of students\u2019 individual paths through particular assignments were qualitatively observed and described.\" \" (2) fine-grain submission rate analysis, where the distribution of submission rates among students was visible.\"\n\nReplace-Text \"The aggregate analysis used data from seven offerings\" \"The analysis used data from seven offerings\"\n\nReplace-Text \"The lack of any higher average submission rates indicates\" \"The lack of significantly higher average submission rates indicates\"\n"}
